# Apply weekly fruit/vegetable price update: swap data rows 2<->4 and 3<->5
# for the date, variety, quality, volume, price min/max/avg, unit and $/kg columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was Golden Nugget / Especial / 44902) becomes the former row 4 data
# (Californiana(o) / Primera / 44505)
$ws.Range("D2").Value = 44505
$ws.Range("K2").Value = "Californiana(o)"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("Q2").Value = "$/bandeja 10 kilos"

# Row 3 (was Golden Nugget / Primera / 44902 / 70-13000) becomes the former
# row 5 data (Golden Nugget / Primera / 44505 / 50-15000)
$ws.Range("D3").Value = 44505
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("S3").Value = 1500

# Row 4 (was Californiana(o) / Primera / 44505) becomes the former row 2 data
# (Golden Nugget / Especial / 44902)
$ws.Range("D4").Value = 44902
$ws.Range("K4").Value = "Golden Nugget"
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 60
$ws.Range("Q4").Value = "$/caja 10 kilos"

# Row 5 (was Golden Nugget / Primera / 44505 / 50-15000) becomes the former
# row 3 data (Golden Nugget / Primera / 44902 / 70-13000)
$ws.Range("D5").Value = 44902
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 13000
$ws.Range("Q5").Value = "$/caja 10 kilos"
$ws.Range("S5").Value = 1300
